# Apply the "Add validation for create test and import question" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix header typo "Answser N" -> "Answer N" (columns K:N on row 5)
$ws.Range("K5").Value = "Answer 1"
$ws.Range("L5").Value = "Answer 2"
$ws.Range("M5").Value = "Answer 3"
$ws.Range("N5").Value = "Answer 4"

# Row 6: update answer 1 / answer 2 sample content
$ws.Range("K6").Value = "fsfsdfsdfsd sdfsdf"
$ws.Range("L6").Value = "sdfsdf"

# Row 9 (Listening / Conversation) - fill in correct answers + validation sample data
$ws.Range("J9").Value = "1,3"
$ws.Range("K9").Value = "sdf"
$ws.Range("L9").Value = "sdfds"
$ws.Range("M9").Value = "sdfsdf"

# Row 10 (Listening / Monologue)
$ws.Range("J10").Value = 3
$ws.Range("K10").Value = "dsfsd"
$ws.Range("L10").Value = "sdfsdf"
$ws.Range("M10").Value = "sdfsdf"

# Row 11 (Writing / Email)
$ws.Range("E11").Value = "source 1"
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = "sdfs"
$ws.Range("L11").Value = "sdfsdf"
$ws.Range("M11").Value = "sf"
$ws.Range("N11").Value = "sf"

# Row 12 (Writing / Essay)
$ws.Range("E12").Value = "source 1"
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = "sdf"
$ws.Range("L12").Value = "sdfsdf"
$ws.Range("M12").Value = "sf"
$ws.Range("N12").Value = "sdf"

# Update the sheet view: scroll so column F is the leftmost visible column,
# and move the active selection to J13.
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("J13").Select()
